$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Cells.Item(87, 25).Value = 1933
$ws.Cells.Item(88, 25).Value = 1934
$ws.Cells.Item(88, 28).Value = 30649
$ws.Cells.Item(89, 25).Value = 1934
$ws.Cells.Item(89, 28).Value = 30669
$ws.Cells.Item(90, 25).Value = 1934
$ws.Cells.Item(90, 28).Value = 30684
$ws.Cells.Item(91, 25).Value = 1935
$ws.Cells.Item(91, 28).Value = 30693
$ws.Cells.Item(92, 23).Value = 93
$ws.Cells.Item(92, 28).Value = 30701

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Cells.Item(88, 25).Value = 151
$ws.Cells.Item(89, 25).Value = 151
$ws.Cells.Item(90, 25).Value = 151
$ws.Cells.Item(91, 25).Value = 151
$ws.Cells.Item(92, 23).Value = 7
$ws.Cells.Item(92, 25).Value = 151

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Cells.Item(5, 25).Value = 7
$ws.Cells.Item(5, 28).Value = 14
$ws.Cells.Item(6, 25).Value = 7
$ws.Cells.Item(6, 28).Value = 16
$ws.Cells.Item(7, 25).Value = 8
$ws.Cells.Item(7, 28).Value = 17
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(8, 28).Value = 23
$ws.Cells.Item(9, 25).Value = 13
$ws.Cells.Item(9, 28).Value = 28
$ws.Cells.Item(10, 25).Value = 14
$ws.Cells.Item(10, 28).Value = 36
$ws.Cells.Item(11, 25).Value = 12
$ws.Cells.Item(11, 28).Value = 40
$ws.Cells.Item(12, 25).Value = 13
$ws.Cells.Item(12, 28).Value = 48
$ws.Cells.Item(13, 25).Value = 13
$ws.Cells.Item(13, 28).Value = 49
$ws.Cells.Item(14, 25).Value = 17
$ws.Cells.Item(14, 28).Value = 62
$ws.Cells.Item(15, 25).Value = 18
$ws.Cells.Item(15, 28).Value = 75
$ws.Cells.Item(16, 25).Value = 21
$ws.Cells.Item(16, 28).Value = 103
$ws.Cells.Item(17, 25).Value = 22
$ws.Cells.Item(17, 28).Value = 124
$ws.Cells.Item(18, 25).Value = 29
$ws.Cells.Item(18, 28).Value = 148
$ws.Cells.Item(19, 25).Value = 31
$ws.Cells.Item(19, 28).Value = 175
$ws.Cells.Item(20, 25).Value = 35
$ws.Cells.Item(20, 28).Value = 277
$ws.Cells.Item(21, 25).Value = 37
$ws.Cells.Item(21, 28).Value = 316
$ws.Cells.Item(22, 25).Value = 43
$ws.Cells.Item(22, 28).Value = 409
$ws.Cells.Item(23, 25).Value = 48
$ws.Cells.Item(23, 28).Value = 501
$ws.Cells.Item(24, 25).Value = 60
$ws.Cells.Item(24, 28).Value = 590
$ws.Cells.Item(25, 25).Value = 65
$ws.Cells.Item(25, 28).Value = 717
$ws.Cells.Item(26, 25).Value = 74
$ws.Cells.Item(26, 28).Value = 820
$ws.Cells.Item(27, 25).Value = 84
$ws.Cells.Item(27, 28).Value = 946
$ws.Cells.Item(28, 25).Value = 93
$ws.Cells.Item(28, 28).Value = 1092
$ws.Cells.Item(29, 25).Value = 105
$ws.Cells.Item(29, 28).Value = 1208
$ws.Cells.Item(30, 25).Value = 119
$ws.Cells.Item(30, 28).Value = 1366
$ws.Cells.Item(31, 25).Value = 120
$ws.Cells.Item(31, 28).Value = 1465
$ws.Cells.Item(32, 25).Value = 130
$ws.Cells.Item(32, 28).Value = 1612
$ws.Cells.Item(33, 25).Value = 139
$ws.Cells.Item(33, 28).Value = 1792
$ws.Cells.Item(34, 25).Value = 144
$ws.Cells.Item(34, 28).Value = 1879
$ws.Cells.Item(35, 25).Value = 154
$ws.Cells.Item(35, 28).Value = 1993
$ws.Cells.Item(36, 25).Value = 155
$ws.Cells.Item(36, 28).Value = 2178
$ws.Cells.Item(37, 25).Value = 154
$ws.Cells.Item(37, 28).Value = 2206
$ws.Cells.Item(38, 25).Value = 148
$ws.Cells.Item(38, 28).Value = 2286
$ws.Cells.Item(39, 25).Value = 148
$ws.Cells.Item(39, 28).Value = 2350
$ws.Cells.Item(40, 25).Value = 148
$ws.Cells.Item(40, 28).Value = 2332
$ws.Cells.Item(41, 25).Value = 151
$ws.Cells.Item(41, 28).Value = 2321
$ws.Cells.Item(42, 25).Value = 148
$ws.Cells.Item(42, 28).Value = 2306
$ws.Cells.Item(43, 25).Value = 135
$ws.Cells.Item(43, 28).Value = 2308
$ws.Cells.Item(44, 25).Value = 132
$ws.Cells.Item(44, 28).Value = 2225
$ws.Cells.Item(45, 25).Value = 125
$ws.Cells.Item(45, 28).Value = 2135
$ws.Cells.Item(46, 25).Value = 121
$ws.Cells.Item(46, 28).Value = 2068
$ws.Cells.Item(47, 25).Value = 119
$ws.Cells.Item(47, 28).Value = 2008
$ws.Cells.Item(48, 25).Value = 116
$ws.Cells.Item(48, 28).Value = 1933
$ws.Cells.Item(49, 25).Value = 114
$ws.Cells.Item(49, 28).Value = 1910
$ws.Cells.Item(50, 25).Value = 110
$ws.Cells.Item(50, 28).Value = 1891
$ws.Cells.Item(51, 25).Value = 104
$ws.Cells.Item(51, 28).Value = 1851
$ws.Cells.Item(52, 25).Value = 98
$ws.Cells.Item(52, 28).Value = 1725
$ws.Cells.Item(53, 25).Value = 96
$ws.Cells.Item(53, 28).Value = 1668
$ws.Cells.Item(54, 25).Value = 89
$ws.Cells.Item(54, 28).Value = 1570
$ws.Cells.Item(55, 25).Value = 88
$ws.Cells.Item(55, 28).Value = 1517
$ws.Cells.Item(56, 25).Value = 86
$ws.Cells.Item(56, 28).Value = 1509
$ws.Cells.Item(57, 25).Value = 83
$ws.Cells.Item(57, 28).Value = 1485
$ws.Cells.Item(58, 25).Value = 80
$ws.Cells.Item(58, 28).Value = 1407
$ws.Cells.Item(59, 25).Value = 81
$ws.Cells.Item(59, 28).Value = 1342
$ws.Cells.Item(60, 25).Value = 77
$ws.Cells.Item(60, 28).Value = 1287
$ws.Cells.Item(61, 25).Value = 75
$ws.Cells.Item(61, 28).Value = 1243
$ws.Cells.Item(62, 25).Value = 77
$ws.Cells.Item(62, 28).Value = 1214
$ws.Cells.Item(63, 25).Value = 79
$ws.Cells.Item(63, 28).Value = 1192
$ws.Cells.Item(64, 25).Value = 70
$ws.Cells.Item(64, 28).Value = 1167
$ws.Cells.Item(65, 25).Value = 66
$ws.Cells.Item(65, 28).Value = 1135
$ws.Cells.Item(66, 25).Value = 66
$ws.Cells.Item(66, 28).Value = 1067
$ws.Cells.Item(67, 25).Value = 66
$ws.Cells.Item(67, 28).Value = 987
$ws.Cells.Item(68, 25).Value = 64
$ws.Cells.Item(68, 28).Value = 932
$ws.Cells.Item(69, 25).Value = 64
$ws.Cells.Item(69, 28).Value = 883
$ws.Cells.Item(70, 25).Value = 65
$ws.Cells.Item(70, 28).Value = 881
$ws.Cells.Item(71, 25).Value = 61
$ws.Cells.Item(71, 28).Value = 857
$ws.Cells.Item(72, 25).Value = 57
$ws.Cells.Item(72, 28).Value = 827
$ws.Cells.Item(73, 25).Value = 54
$ws.Cells.Item(73, 28).Value = 798
$ws.Cells.Item(74, 25).Value = 49
$ws.Cells.Item(74, 28).Value = 747
$ws.Cells.Item(75, 25).Value = 47
$ws.Cells.Item(75, 28).Value = 700
$ws.Cells.Item(76, 25).Value = 46
$ws.Cells.Item(76, 28).Value = 669
$ws.Cells.Item(77, 25).Value = 47
$ws.Cells.Item(77, 28).Value = 658
$ws.Cells.Item(78, 25).Value = 47
$ws.Cells.Item(78, 28).Value = 649
$ws.Cells.Item(79, 25).Value = 45
$ws.Cells.Item(79, 28).Value = 646
$ws.Cells.Item(80, 25).Value = 42
$ws.Cells.Item(80, 28).Value = 596
$ws.Cells.Item(81, 25).Value = 40
$ws.Cells.Item(81, 28).Value = 560
$ws.Cells.Item(82, 25).Value = 38
$ws.Cells.Item(82, 28).Value = 524
$ws.Cells.Item(83, 25).Value = 38
$ws.Cells.Item(83, 28).Value = 510
$ws.Cells.Item(84, 25).Value = 38
$ws.Cells.Item(84, 28).Value = 512
$ws.Cells.Item(85, 25).Value = 36
$ws.Cells.Item(85, 28).Value = 494
$ws.Cells.Item(86, 25).Value = 35
$ws.Cells.Item(86, 28).Value = 478
$ws.Cells.Item(87, 25).Value = 33
$ws.Cells.Item(87, 28).Value = 463
$ws.Cells.Item(88, 25).Value = 33
$ws.Cells.Item(88, 28).Value = 440
$ws.Cells.Item(89, 25).Value = 32
$ws.Cells.Item(89, 28).Value = 426
$ws.Cells.Item(90, 25).Value = 32
$ws.Cells.Item(90, 28).Value = 413
$ws.Cells.Item(91, 25).Value = 32
$ws.Cells.Item(91, 28).Value = 419
$ws.Cells.Item(92, 23).Value = 1
$ws.Cells.Item(92, 25).Value = 32
$ws.Cells.Item(92, 28).Value = 407

$ws = $wb.Worksheets.Item("ICU")
$ws.Cells.Item(87, 25).Value = 5
$ws.Cells.Item(87, 28).Value = 50
$ws.Cells.Item(88, 25).Value = 5
$ws.Cells.Item(88, 28).Value = 48
$ws.Cells.Item(89, 25).Value = 4
$ws.Cells.Item(89, 28).Value = 47
$ws.Cells.Item(90, 25).Value = 4
$ws.Cells.Item(90, 28).Value = 47
$ws.Cells.Item(91, 25).Value = 4
$ws.Cells.Item(91, 28).Value = 44
$ws.Cells.Item(92, 25).Value = 4
$ws.Cells.Item(92, 28).Value = 42

$ws = $wb.Worksheets.Item("Ventilated")
$ws.Cells.Item(88, 25).Value = 4
$ws.Cells.Item(89, 25).Value = 3
$ws.Cells.Item(89, 28).Value = 16
$ws.Cells.Item(90, 25).Value = 3
$ws.Cells.Item(90, 28).Value = 15
$ws.Cells.Item(91, 25).Value = 3
$ws.Cells.Item(91, 28).Value = 14
$ws.Cells.Item(92, 25).Value = 3
$ws.Cells.Item(92, 28).Value = 12

$ws = $wb.Worksheets.Item("Released")
$ws.Cells.Item(92, 23).Value = 84
